$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Helper: add a bookmark at a single position and immediately delete
# it again.  Adding (then deleting) a bookmark forces the engine to
# keep the text on either side of that position in separate <w:r>
# runs without touching run formatting - exactly mirroring the way
# Word keeps runs that were typed/edited separately from merging
# back together.
# -----------------------------------------------------------------
function Split-At([int]$pos, [string]$bmName) {
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add($bmName, $bmRange) | Out-Null
    $d.Bookmarks($bmName).Delete()
}

# ===================================================================
# Change 1 - "Procedure" section: insert " Euclidian" between
# "3D" and " distance" so the sentence reads "...the 3D Euclidian
# distance...", as three separate runs.
# ===================================================================
$find1 = $d.Content.Find
$found1 = $find1.Execute("we used the 3D", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'we used the 3D' anchor text"
}
$afterThreeD = $find1.Parent.End

$insRange = $d.Range($afterThreeD, $afterThreeD)
$insRange.InsertAfter(" Euclidian")
$afterEuclidian = $insRange.End

Split-At $afterThreeD "ZZZ_SPLIT_1"
Split-At $afterEuclidian "ZZZ_SPLIT_2"

# ===================================================================
# Change 2 - "Future Ideas" section: replace "In addition" with
# "Finally" and move the hidden _GoBack bookmark so that it now sits
# right after "Finally" (mirroring a user selecting "In addition"
# and typing "Finally" over it).
# ===================================================================
$find2 = $d.Content.Find
$found2 = $find2.Execute("One way that we could improve our algorithm", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find target paragraph anchor text"
}
$paraRange = $find2.Parent.Paragraphs(1).Range

$scopeRange = $d.Range($paraRange.Start, $paraRange.End)
$find3 = $scopeRange.Find
$found3 = $find3.Execute("In addition", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not find 'In addition' inside target paragraph"
}
$targetRange = $find3.Parent
$targetRange.Text = "Finally"
$finallyStart = $targetRange.Start
$finallyEnd = $targetRange.End

# Locate ", we could use " right after the freshly typed "Finally"
$scope2 = $d.Range($finallyEnd, $d.Content.End)
$find4 = $scope2.Find
$found4 = $find4.Execute(", we could use ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found4) {
    throw "Could not find ', we could use ' after replacement"
}
$commaStart = $find4.Parent.Start
$afterComma = $commaStart + 2   # length of ", "

# Locate the end of "we could use " (start of "physical proximity...")
$scope3 = $d.Range($afterComma, $d.Content.End)
$find5 = $scope3.Find
$found5 = $find5.Execute("we could use ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found5) {
    throw "Could not find 'we could use ' after replacement"
}
$weCouldUseEnd = $find5.Parent.End

# Split the paragraph into the five runs the target structure needs,
# and drop the (now-relocated) hidden _GoBack bookmark right after
# "Finally".
Split-At $finallyStart "ZZZ_SPLIT_3"
$d.Bookmarks.Add("_GoBack", $d.Range($finallyEnd, $finallyEnd)) | Out-Null
Split-At $afterComma "ZZZ_SPLIT_4"
Split-At $weCouldUseEnd "ZZZ_SPLIT_5"

Write-Host "Edits applied successfully."
